$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 517.8
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 599.75
$ws.Range("K2").Value = 190
$ws.Range("L2").Value = 599.75
$ws.Range("M2").Value = -77
$ws.Range("N2").Value = -825.75
$ws.Range("H9").Value = 365.30768
$ws.Range("I9").Value = 256.125
$ws.Range("K9").Value = 256.125
$ws.Range("M9").Value = -87.125
$ws.Range("H32").Value = 11114163
$ws.Range("J32").Value = 14289481
$ws.Range("L32").Value = 14289481
$ws.Range("N32").Value = -14290133
$ws.Range("H62").Value = 3484.5
$ws.Range("J62").Value = 3989.6667
$ws.Range("L62").Value = 3989.6667
$ws.Range("N62").Value = -5237.6667
$ws.Range("H65").Value = 3484.5
$ws.Range("J65").Value = 3989.6667
$ws.Range("L65").Value = 19948.3335
$ws.Range("N65").Value = -26188.3335
$ws.Range("H86").Value = 4223.8184
$ws.Range("I86").Value = 2994.5715
$ws.Range("K86").Value = 2994.5715
$ws.Range("M86").Value = -1871.5715
$ws.Range("H89").Value = 4223.8184
$ws.Range("I89").Value = 2994.5715
$ws.Range("K89").Value = 14972.8575
$ws.Range("M89").Value = -9356.8575
$ws.Range("H116").Value = 10999.728
$ws.Range("J116").Value = 11332.833
$ws.Range("L116").Value = 11332.833
$ws.Range("N116").Value = -18216.833
$ws.Range("H132").Value = 14504.619
$ws.Range("I132").Value = 1900.9667
$ws.Range("J132").Value = 46013.75
$ws.Range("K132").Value = 5702.9001
$ws.Range("L132").Value = 138041.25
$ws.Range("M132").Value = -3172.9001
$ws.Range("N132").Value = -143101.25
$ws.Range("H137").Value = 2868.4285
$ws.Range("I137").Value = 2731.1
$ws.Range("K137").Value = 8193.299999999999
$ws.Range("M137").Value = -5643.299999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3090
$ws.Range("I61").Value = 2709
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 2709
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -2497
$ws.Range("N61").Value = -5419
$ws.Range("H136").Value = 3090
$ws.Range("I136").Value = 2709
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 8127
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -5577
$ws.Range("N136").Value = -20085

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20835258
$ws.Range("I86").Value = 38463476
$ws.Range("J86").Value = 1907.2727
$ws.Range("K86").Value = 38463476
$ws.Range("L86").Value = 1907.2727
$ws.Range("M86").Value = -38462353
$ws.Range("N86").Value = -4153.2727
$ws.Range("H89").Value = 20835258
$ws.Range("I89").Value = 38463476
$ws.Range("J89").Value = 1907.2727
$ws.Range("K89").Value = 192317380
$ws.Range("L89").Value = 9536.363499999999
$ws.Range("M89").Value = -192311764
$ws.Range("N89").Value = -20768.3635

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3544.0952
$ws.Range("I99").Value = 3707.8
$ws.Range("K99").Value = 3707.8
$ws.Range("M99").Value = -2209.8
$ws.Range("H107").Value = 8929624
$ws.Range("I107").Value = 12987925
$ws.Range("J107").Value = 1361.4
$ws.Range("K107").Value = 12987925
$ws.Range("L107").Value = 1361.4
$ws.Range("M107").Value = -12986005
$ws.Range("N107").Value = -5201.4
$ws.Range("H122").Value = 2746.4375
$ws.Range("I122").Value = 3150
$ws.Range("J122").Value = 1976
$ws.Range("K122").Value = 9450
$ws.Range("L122").Value = 5928
$ws.Range("M122").Value = -7000
$ws.Range("N122").Value = -10828
$ws.Range("H126").Value = 3544.0952
$ws.Range("I126").Value = 3707.8
$ws.Range("K126").Value = 11123.4
$ws.Range("M126").Value = -8653.400000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2408.261
$ws.Range("J55").Value = 4500
$ws.Range("L55").Value = 13500
$ws.Range("N55").Value = -13854
$ws.Range("H80").Value = 6175
$ws.Range("I80").Value = 7500
$ws.Range("J80").Value = 5607.143
$ws.Range("K80").Value = 22500
$ws.Range("L80").Value = 16821.429
$ws.Range("M80").Value = -21564
$ws.Range("N80").Value = -18693.429
$ws.Range("H83").Value = 6175
$ws.Range("I83").Value = 7500
$ws.Range("J83").Value = 5607.143
$ws.Range("K83").Value = 67500
$ws.Range("L83").Value = 50464.287
$ws.Range("M83").Value = -62820
$ws.Range("N83").Value = -59824.287
$ws.Range("H92").Value = 457.66666
$ws.Range("J92").Value = 457.66666
$ws.Range("L92").Value = 1372.99998
$ws.Range("N92").Value = -3868.99998
$ws.Range("H131").Value = 3153.8333
$ws.Range("I131").Value = 1603
$ws.Range("J131").Value = 5590.857
$ws.Range("K131").Value = 4809
$ws.Range("L131").Value = 16772.571
$ws.Range("M131").Value = 231
$ws.Range("N131").Value = -26852.571
$ws.Range("H132").Value = 2246.0667
$ws.Range("I132").Value = 2731.1667
$ws.Range("J132").Value = 1922.6666
$ws.Range("K132").Value = 24580.5003
$ws.Range("L132").Value = 17303.9994
$ws.Range("M132").Value = -22050.5003
$ws.Range("N132").Value = -22363.9994

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 89216.62
$ws.Range("I70").Value = 161923.42
$ws.Range("J70").Value = 4392
$ws.Range("K70").Value = 161923.42
$ws.Range("L70").Value = 4392
$ws.Range("M70").Value = -161653.42
$ws.Range("N70").Value = -4932
$ws.Range("H73").Value = 89216.62
$ws.Range("I73").Value = 161923.42
$ws.Range("J73").Value = 4392
$ws.Range("K73").Value = 161923.42
$ws.Range("L73").Value = 4392
$ws.Range("M73").Value = -160987.42
$ws.Range("N73").Value = -6264
$ws.Range("H102").Value = 5126.1577
$ws.Range("I102").Value = 4524.8125
$ws.Range("K102").Value = 4524.8125
$ws.Range("M102").Value = -2902.8125
$ws.Range("H107").Value = 600.12
$ws.Range("I107").Value = 504.89474
$ws.Range("J107").Value = 901.6667
$ws.Range("K107").Value = 504.89474
$ws.Range("L107").Value = 901.6667
$ws.Range("M107").Value = 1415.10526
$ws.Range("N107").Value = -4741.6667
$ws.Range("H126").Value = 705
$ws.Range("I126").Value = 705
$ws.Range("K126").Value = 2115
$ws.Range("M126").Value = 355

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7417
$ws.Range("J7").Value = 8999.666999999999
$ws.Range("L7").Value = 8999.666999999999
$ws.Range("N7").Value = -9223.666999999999
$ws.Range("H40").Value = 2617.625
$ws.Range("I40").Value = 2098.7144
$ws.Range("K40").Value = 2098.7144
$ws.Range("M40").Value = -1962.7144
$ws.Range("H61").Value = 2153.2727
$ws.Range("I61").Value = 2148.375
$ws.Range("J61").Value = 2166.3333
$ws.Range("K61").Value = 2148.375
$ws.Range("L61").Value = 2166.3333
$ws.Range("M61").Value = -1946.375
$ws.Range("N61").Value = -2570.3333
$ws.Range("H113").Value = 2153.2727
$ws.Range("I113").Value = 2148.375
$ws.Range("J113").Value = 2166.3333
$ws.Range("K113").Value = 2148.375
$ws.Range("L113").Value = 2166.3333
$ws.Range("M113").Value = 21.625
$ws.Range("N113").Value = -6506.3333
$ws.Range("H126").Value = 7417
$ws.Range("J126").Value = 8999.666999999999
$ws.Range("L126").Value = 26999.001
$ws.Range("N126").Value = -31939.001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10758.625
$ws.Range("J41").Value = 11152.714
$ws.Range("L41").Value = 11152.714
$ws.Range("N41").Value = -11932.714
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51144
$ws.Range("H107").Value = 1434.1666
$ws.Range("I107").Value = 2491
$ws.Range("J107").Value = 905.75
$ws.Range("K107").Value = 7473
$ws.Range("L107").Value = 2717.25
$ws.Range("M107").Value = -5553
$ws.Range("N107").Value = -6557.25
